# Applies the "Updated symbol list" crypto price/volume refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price), E (Volume 1h %) and G (Hora) hold numeric-looking text
# (e.g. "256.08", "-0.09%", "14"). Force them to Text format before writing
# so Excel keeps the literal strings instead of coercing to numbers.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$changes = @(
    @('D2', '256.08'),
    @('E2', '-0.09%'),
    @('G2', '14'),
    @('D3', '26.56'),
    @('E3', '-1.94%'),
    @('G3', '14'),
    @('D4', '4.647'),
    @('E4', '0.25%'),
    @('G4', '14'),
    @('D5', '0.05921'),
    @('E5', '0.40%'),
    @('G5', '14'),
    @('E6', '-0.48%'),
    @('G6', '14'),
    @('D7', '0.8550'),
    @('E7', '-1.51%'),
    @('G7', '14'),
    @('D8', '0.9121'),
    @('E8', '-3.44%'),
    @('G8', '14'),
    @('B9', 'One'),
    @('C9', 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'),
    @('D9', '0.01031'),
    @('E9', '1,612.77%'),
    @('G9', '14'),
    @('B10', 'WazirX'),
    @('C10', 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'),
    @('D10', '0.1377'),
    @('E10', '-1.95%'),
    @('G10', '14'),
    @('B11', 'LiechtensteinCryptoassetsExchange'),
    @('C11', 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'),
    @('D11', '0.04218'),
    @('E11', '13.19%'),
    @('G11', '14'),
    @('B12', 'MandalaExchangeToken'),
    @('C12', 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'),
    @('D12', '0.07007'),
    @('E12', '-0.98%'),
    @('G12', '14'),
    @('B13', 'BitrueCoin'),
    @('C13', 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'),
    @('D13', '0.03037'),
    @('E13', '-4.94%'),
    @('G13', '14'),
    @('B14', 'BitMartToken'),
    @('C14', 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'),
    @('D14', '0.09107'),
    @('E14', '-1.61%'),
    @('G14', '14'),
    @('B15', 'BitForexToken'),
    @('C15', 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'),
    @('D15', '0.001529'),
    @('E15', '-1.30%'),
    @('G15', '14'),
    @('B16', 'TigerCash'),
    @('C16', 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'),
    @('D16', '0.006070'),
    @('E16', '-0.11%'),
    @('G16', '14'),
    @('B17', 'LEO'),
    @('C17', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'),
    @('D17', '3.472'),
    @('E17', '-1.12%'),
    @('G17', '14'),
    @('B18', 'GateToken'),
    @('C18', 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'),
    @('D18', '3.141'),
    @('E18', '-1.57%'),
    @('G18', '14'),
    @('B19', 'BTSEToken'),
    @('C19', 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'),
    @('D19', '2.150'),
    @('E19', '-2.37%'),
    @('G19', '14'),
    @('D20', '0.3086'),
    @('E20', '-0.62%'),
    @('G20', '14'),
    @('E21', '0.08%'),
    @('G21', '14'),
    @('D22', '3.871'),
    @('E22', '0.47%'),
    @('G22', '14'),
    @('D23', '0.04222'),
    @('E23', '-0.29%'),
    @('G23', '14'),
    @('D24', '0.001217'),
    @('E24', '-0.23%'),
    @('G24', '14'),
    @('D25', '0.004648'),
    @('E25', '8.55%'),
    @('G25', '14'),
    @('D26', '0.0001200'),
    @('E26', '-0.02%'),
    @('G26', '14'),
    @('D27', '0.0001715'),
    @('E27', '14.20%'),
    @('G27', '14'),
    @('G28', '14'),
    @('G29', '14'),
    @('G30', '14'),
    @('G31', '14'),
    @('G32', '14'),
    @('G33', '14'),
    @('G34', '14'),
    @('G35', '14'),
    @('G36', '14'),
    @('G37', '14'),
    @('G38', '14'),
    @('G39', '14'),
    @('D40', '0.03790'),
    @('E40', '-0.53%'),
    @('G40', '14'),
    @('D41', '0.006179'),
    @('E41', '0.31%'),
    @('G41', '14'),
    @('D42', '0.1098'),
    @('E42', '-0.10%'),
    @('G42', '14'),
    @('D43', '0.002311'),
    @('E43', '21.62%'),
    @('G43', '14'),
    @('D44', '0.01449'),
    @('E44', '28.29%'),
    @('G44', '14'),
    @('D45', '0.00005131'),
    @('E45', '-6.71%'),
    @('G45', '14'),
    @('D46', '0.00000000750'),
    @('E46', '-0.02%'),
    @('G46', '14'),
    @('D47', '0.04999'),
    @('E47', '-35.75%'),
    @('G47', '14'),
    @('E48', '10,463.61%'),
    @('G48', '14'),
    @('D49', '0.00002100'),
    @('E49', '-0.02%'),
    @('G49', '14'),
    @('D50', '0.0002000'),
    @('E50', '-0.02%'),
    @('G50', '14'),
    @('G51', '14')
)

foreach ($pair in $changes) {
    $ws.Range($pair[0]).Value = $pair[1]
}
